$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.309.72'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.867.75'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.34'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4806'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2797'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06492'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.859.80'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07446'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.44'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.061'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.78'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6524'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.287.67'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.23'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007575'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.107.86'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.12%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.273'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '218.65'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +13.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.142'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.287'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.61'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.38'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.965'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.452'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09353'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.293'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.006'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05035'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.199'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +9.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7444'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.711'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01817'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.614'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.072'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9069'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.94%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.913'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.37%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.47'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.003'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4249'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.361'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1283'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '63.66'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.473'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.902'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.53'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.71%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05629'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.95%  '
